$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246; this shifts existing rows 246:287 down to 247:288
$ws.Rows(246).Insert()

# Populate the new row 246 with a new price record (New Hall variety)
$ws.Cells.Item(246, 1).Value = 11
$ws.Cells.Item(246, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(246, 3).Value = "Bíobío"
$ws.Cells.Item(246, 4).Value = 44722
$ws.Cells.Item(246, 5).Value = 8
$ws.Cells.Item(246, 6).Value = "Fruta"
$ws.Cells.Item(246, 7).Value = 100102
$ws.Cells.Item(246, 8).Value = "Cítricos"
$ws.Cells.Item(246, 9).Value = 100102005
$ws.Cells.Item(246, 10).Value = "Naranja"
$ws.Cells.Item(246, 11).Value = "New Hall"
$ws.Cells.Item(246, 12).Value = "Primera"
$ws.Cells.Item(246, 13).Value = 400
$ws.Cells.Item(246, 14).Value = 8000
$ws.Cells.Item(246, 15).Value = 9000
$ws.Cells.Item(246, 16).Value = 8500
$ws.Cells.Item(246, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(246, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(246, 19).Value = 567
$ws.Cells.Item(246, 20).Value = 15
